$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Orden" column (I) -------------------------------------------------
# Copy the header formatting from H1 (bold/centered/bordered) onto I1, then
# set its text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "Orden"

# Fill in the "Orden" value for the existing 18 products (rows 2-19): 8..25
$orden = 8
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 9).Value = $orden
    $orden++
}

# --- New stock / products (ids 19-25 -> sheet rows 20-26) ------------------
$newRows = @(
    @{ Id=19; Nombre="Pijama 19"; Valor=10000; Stock=5; Img="productos/producto-1.webp"; Desc=$true;  Video="/productos/producto-1.mp4"; Orden=1 },
    @{ Id=20; Nombre="Pijama 20"; Valor=10000; Stock=5; Img="productos/producto-2.webp"; Desc=$false; Video="/Producto1.mp4";              Orden=4 },
    @{ Id=21; Nombre="Pijama 21"; Valor=6000;  Stock=5; Img="productos/producto-3.webp"; Desc=$false; Video="/productos/producto-2.mp4";   Orden=3 },
    @{ Id=22; Nombre="Pijama 22"; Valor=3000;  Stock=5; Img="productos/producto-4.webp"; Desc=$false; Video="/productos/producto-2.mp4";   Orden=2 },
    @{ Id=23; Nombre="Pijama 23"; Valor=8000;  Stock=5; Img="productos/producto-5.webp"; Desc=$false; Video="/Producto1.mp4";              Orden=5 },
    @{ Id=24; Nombre="Pijama 24"; Valor=15000; Stock=1; Img="productos/producto-6.webp"; Desc=$false; Video="/Producto1.mp4";              Orden=6 },
    @{ Id=25; Nombre="Pijama 25"; Valor=10000; Stock=3; Img="productos/producto-7.webp"; Desc=$false; Video="/Producto1.mp4";              Orden=7 }
)

$row = 20
foreach ($p in $newRows) {
    $ws.Cells.Item($row, 1).Value = $p.Id
    $ws.Cells.Item($row, 2).Value = $p.Nombre
    $ws.Cells.Item($row, 3).Value = "Diseñados en colombia, piel de durazno."
    $ws.Cells.Item($row, 4).Value = $p.Valor
    $ws.Cells.Item($row, 5).Value = $p.Stock
    $ws.Cells.Item($row, 6).Value = $p.Img
    $ws.Cells.Item($row, 7).Value = $p.Desc
    $ws.Cells.Item($row, 8).Value = $p.Video
    $ws.Cells.Item($row, 9).Value = $p.Orden
    $row++
}

# --- Column widths (mirrors the real workbook's "AutoFit selected columns") -
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666   # A - IdProducto
$ws.Columns.Item(2).ColumnWidth = 15.0                # B - NombreProducto
$ws.Columns.Item(3).ColumnWidth = 33.333333333333336  # C - Descripcion
$ws.Columns.Item(6).ColumnWidth = 23.0                # F - ImageUrl
$ws.Columns.Item(7).ColumnWidth = 12.5                # G - ConDescuento
$ws.Columns.Item(8).ColumnWidth = 22.833333333333332  # H - VideoUrl

# --- Selection matches where the author ended up after the edit -----------
$ws.Range("H26").Select() | Out-Null
